$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-6 get rotated: data that was in row 3 moves to row 6,
# row 4 -> row 5, row 5 -> row 3, row 6 -> row 4.
# Only columns D, L, M, N, O, P, Q, R, S, T change (others stay identical
# across these rows already).

$rowData = @{
    3 = @{ D = 44316; L = "Primera"; M = 60;  N = 17500; O = 18000; P = 17750; Q = "`$/caja 16 kilos granel";       R = "Región de O'Higgins"; S = 1109; T = 16 }
    4 = @{ D = 44316; L = "Segunda"; M = 40;  N = 16000; O = 16000; P = 16000; Q = "`$/caja 16 kilos granel";       R = "Región de O'Higgins"; S = 1000; T = 16 }
    5 = @{ D = 44344; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; Q = "`$/caja 18 kilos granel";       R = "Provincia de Curicó";  S = 750;  T = 18 }
    6 = @{ D = 44334; L = "Primera"; M = 120; N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada";    R = "Región de O'Higgins"; S = 1042; T = 12 }
}

foreach ($r in 3..6) {
    $vals = $rowData[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
